$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(4).Insert()
$w = $ws.Columns.Item(5).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $w
